# Commit: "update database and change read_price algorithm"
#
# The yearly period columns (E:I) roll forward by one fiscal year:
#   old E..I = FY1396/12, FY1397/12, FY1398/12, FY1399/12, FY1400/12
#   new E..I = FY1397/12, FY1398/12, FY1399/12, FY1400/12, FY1401/12
# i.e. every data column shifts one slot to the left (oldest year dropped)
# and a brand-new FY1401/12 column of data is appended on the right (I).
#
# This applies both to the two header rows (8 and 24, "دوازده ماهه منتهی
# به ..." labels) and to every numeric data row beneath each header.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header labels (row 8 and row 24) ----------------------------------
$headerLabels = @(
    "دوازده ماهه منتهی به 1397/12",
    "دوازده ماهه منتهی به 1398/12",
    "دوازده ماهه منتهی به 1399/12",
    "دوازده ماهه منتهی به 1400/12",
    "دوازده ماهه منتهی به 1401/12"
)
$cols = @("E", "F", "G", "H", "I")

for ($i = 0; $i -lt 5; $i++) {
    $ws.Range($cols[$i] + "8").Value = $headerLabels[$i]
    $ws.Range($cols[$i] + "24").Value = $headerLabels[$i]
}

# ---- Numeric data rows: shift E<-F<-G<-H<-I, new value lands in I ------
# New rightmost-column (I) values introduced by this edit, per row.
$newRightValues = @{
    10 = 0
    11 = 0
    12 = 457337
    13 = 0
    14 = 0
    15 = 0
    16 = 29844
    17 = 893395
    18 = 0
    19 = 1557227
    20 = 2937803
    26 = 743
    27 = 935
}

foreach ($r in $newRightValues.Keys) {
    $oldVals = @(
        $ws.Range("E$r").Value(),
        $ws.Range("F$r").Value(),
        $ws.Range("G$r").Value(),
        $ws.Range("H$r").Value(),
        $ws.Range("I$r").Value()
    )
    # shift left by one
    $ws.Range("E$r").Value = $oldVals[1]
    $ws.Range("F$r").Value = $oldVals[2]
    $ws.Range("G$r").Value = $oldVals[3]
    $ws.Range("H$r").Value = $oldVals[4]
    # append the brand-new year's figure on the right
    $ws.Range("I$r").Value = $newRightValues[$r]
}
